$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> [Price(D), Volume(E)] updates. Empty string for D means "no change".
$updates = @{
    2  = @("22.155.49", "  -1.42%  ")
    3  = @("1.560.93", "  -0.96%  ")
    4  = @("", "  -0.04%  ")
    5  = @("", "  +0.00%  ")
    6  = @("290.11", "  +0.46%  ")
    7  = @("0.3801", "  +3.17%  ")
    8  = @("0.3283", "  -1.51%  ")
    9  = @("43.74", "  -8.99%  ")
    10 = @("1.138", "  -0.79%  ")
    11 = @("0.07362", "  -2.55%  ")
    12 = @("1.001", "  -0.06%  ")
    13 = @("19.95", "  -4.16%  ")
    14 = @("5.825", "  -2.31%  ")
    15 = @("6.869", "  -1.13%  ")
    16 = @("1.563.47", "  -0.73%  ")
    17 = @("0.00001094", "  -2.60%  ")
    18 = @("0.06631", "  -1.46%  ")
    19 = @("85.46", "  -2.74%  ")
    20 = @("6.451", "  +0.83%  ")
    22 = @("", "  -2.89%  ")
    23 = @("11.74", "  -2.21%  ")
    24 = @("22.163.96", "  -1.33%  ")
    25 = @("2.266", "  -5.14%  ")
    26 = @("2.535", "  -4.17%  ")
    27 = @("", "  -0.13%  ")
    28 = @("19.07", "  -3.13%  ")
    29 = @("4.869", "  -2.60%  ")
    30 = @("1.738.83", "  -0.73%  ")
    31 = @("121.38", "  -3.20%  ")
    32 = @("1.119", "  +1.98%  ")
    33 = @("6.036", "  -1.83%  ")
    34 = @("1.884", "  -5.56%  ")
    35 = @("9.345", "  -5.48%  ")
    36 = @("0.08191", "  -2.19%  ")
    37 = @("5.285", "  -1.76%  ")
    38 = @("0.02305", "  -6.88%  ")
    39 = @("0.06218", "  -3.00%  ")
    40 = @("0.2139", "  -5.04%  ")
    41 = @("1.235", "  -4.48%  ")
    42 = @("11.07", "  -3.57%  ")
    43 = @("", "  +0.00%  ")
    44 = @("0.5981", "  -5.18%  ")
    45 = @("13.67", "  -3.51%  ")
    46 = @("3.757", "  -0.81%  ")
    47 = @("0.5789", "  -5.72%  ")
    48 = @("1.986", "  -3.98%  ")
    49 = @("120.95", "  -3.82%  ")
    50 = @("1.170", "  -3.71%  ")
    51 = @("0.06995", "  -3.34%  ")
}

foreach ($row in $updates.Keys) {
    $pair = $updates[$row]
    $price = $pair[0]
    $volume = $pair[1]

    if ($price -ne "") {
        # These "Price" cells are stored as plain text in the sheet (several
        # contain two '.' separators, e.g. "22.155.49", so they can never be
        # real numbers). Force the cell to Text format first so Excel does
        # not auto-convert numeric-looking strings like "290.11" into a
        # floating point number, then clear the formatting again so the
        # cell keeps its original (default) style -- only the text value
        # itself should change.
        $cell = $ws.Range("D$row")
        $cell.NumberFormat = "@"
        $cell.Value = $price
        $cell.ClearFormats()
    }
    $ws.Range("E$row").Value = $volume
}
